$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing card "Look Into the Future" -> "Foresight" (row 72, column A)
$ws.Cells.Item(72, 1).Value = "Foresight"

# Append new "Generic" cards in rows 75-78.
# Write the ImageFile (col C) codes first, then the Names (col A), then the
# shared Kingdom/Type/Rarity values, so new shared-string entries land in the
# same order the original author's save produced.
$ws.Cells.Item(75, 3).Value = "S223"
$ws.Cells.Item(76, 3).Value = "S224"
$ws.Cells.Item(77, 3).Value = "S225"
$ws.Cells.Item(78, 3).Value = "S226"

$ws.Cells.Item(75, 1).Value = "Chaos Wolf"
$ws.Cells.Item(76, 1).Value = "Flaming Eagle"
$ws.Cells.Item(77, 1).Value = "Ruins of the Night"
$ws.Cells.Item(78, 1).Value = "Shrine of Everlasting Day"

$ws.Cells.Item(75, 2).Value = "Generic"
$ws.Cells.Item(76, 2).Value = "Generic"
$ws.Cells.Item(77, 2).Value = "Generic"
$ws.Cells.Item(78, 2).Value = "Generic"

$ws.Cells.Item(75, 4).Value = "Creature"
$ws.Cells.Item(76, 4).Value = "Creature"
$ws.Cells.Item(77, 4).Value = "Structure"
$ws.Cells.Item(78, 4).Value = "Structure"

$ws.Cells.Item(75, 5).Value = "Common"
$ws.Cells.Item(76, 5).Value = "Common"
$ws.Cells.Item(77, 5).Value = "Ultra Rare"
$ws.Cells.Item(78, 5).Value = "Ultra Rare"

# Move the selection/view to just past the last row, matching the author's
# saved cursor position after adding the new cards.
$ws.Range("A79").Select()
